# Word COM-interop script implementing the diff:
#  1. Insert a new paragraph after the title ("Play Bomber Fruit Free - Read Our
#     Slot Review") containing a bold "Meta description" run followed by a
#     normal ": Read our review of Bomber Fruit..." run.
#  2. Remove the duplicate bold "Play Bomber Fruit Free - Read Our Slot Review"
#     paragraph that used to sit right before the final (italic) paragraph.
#  3. Replace the text of the final italic paragraph with the DALLE image
#     prompt, keeping its italic formatting intact.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Insert the "Meta description" paragraph right after the first paragraph.
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()

$metaPara = $d.Paragraphs(2)
# Strip the inherited Heading1 style so the new paragraph is plain body text.
$metaPara.Range.Style = "Normal"

$metaLabel = "Meta description"
$metaRest  = ": Read our review of Bomber Fruit, a fruit-themed online slot game featuring wild symbols and an above average RTP. Play for free!"

$insPoint = $metaPara.Range
$insPoint.Collapse(1)
$insPoint.InsertAfter($metaLabel + $metaRest)

# Bold only the "Meta description" label portion.
$metaPara = $d.Paragraphs(2)
$boldRange = $d.Range($metaPara.Range.Start, $metaPara.Range.Start + $metaLabel.Length)
$boldRange.Bold = 1

# ---------------------------------------------------------------------------
# 2 & 3. Near the end of the document: drop the duplicate bold title
#         paragraph and rewrite the italic paragraph's text.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$dupTitlePara = $d.Paragraphs($count - 1)
$dupTitlePara.Range.Delete()

$count = $d.Paragraphs.Count
$italicPara = $d.Paragraphs($count)
$italicRange = $italicPara.Range
# Exclude the trailing paragraph mark from the replaced text.
$textOnly = $d.Range($italicRange.Start, $italicRange.End - 1)
$textOnly.Text = "Prompt: DALLE, create a cartoon-style image featuring a happy Maya warrior with glasses for the game `"Bomber Fruit`". The image should include the Maya warrior holding a bomb and a basket of fruits in his other hand. The background should be colorful and incorporate some of the fruits from the game such as watermelon, lemon, and cherry. Make sure the image is eye-catching and playful, while still capturing the essence of the game."
